$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "Kernelized SVM" entry (row 8 in the table) right after the
# existing "Support Vector Machine" (SVM) entry. The new block re-uses the
# exact same visual style as the SVM block above it (rows 18-19), so the
# cheapest and most faithful way to reproduce it is to copy that block's
# formatting down onto the two new rows, then overwrite the text.
# ---------------------------------------------------------------------------

# 1) Clone the look (fonts, colors, fills, alignment) of the SVM block
#    (B18:G19) onto the two brand-new rows (22:23).
$ws.Range("B18:G19").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Match the row heights used by every other data row in the table.
$ws.Rows.Item(22).RowHeight = 16.5
$ws.Rows.Item(23).RowHeight = 16.5

# 3) Fill in the new content.
$ws.Range("B22").Value = 8
$ws.Range("C22").Value = "Kernelized SVM"
$ws.Range("D22").Value = "from sklearn.svm import LinearSVC"
$ws.Range("E22").Value = "C & gamma"
$ws.Range("F22").Value = "Better C and gamma means more accurate model"
$ws.Range("G22").Value = "Regression / Classification"

$ws.Range("C23").Value = "(SVM)"
$ws.Range("D23").Value = "svm = SVC(kernel='rbf', C=10, gamma=0.1).fit(X, y)"

# The second sub-row's "(SVM)" label in C23 didn't inherit the colored
# C-column label style from the paste (that style only lives on the first
# sub-row in each block); match the plain bold-label look used elsewhere in
# the sheet for this kind of secondary caption.
$ws.Range("C23").Font.Name = "Arial"
$ws.Range("C23").Font.Bold = $true
$ws.Range("C23").Font.Color = 0
$ws.Range("C23").Interior.Pattern = -4142

# 4) Move the active selection to G22, matching where editing finished.
[void]$ws.Range("G22").Select()
